$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("newsheet")

# Row 1: keep "abc" in A1, set class/lesson info in B1/C1
$ws.Range("B1").Value = "'2/1"
$ws.Range("C1").Value = "phương nguyên"

# Row 2 and Row 3: clear entirely (no more data there)
$ws.Range("A2:C3").Clear()

# Row 4: new class/lesson entry, with "abc" now also present in A4
$ws.Range("A4").Value = "abc"
$ws.Range("B4").Value = "'3/4"
$ws.Range("C4").Value = "phương nguyên"

# Row 6: clear entirely
$ws.Range("A6:C6").Clear()

# Row 7: new class/lesson entry added next to existing "abc"
$ws.Range("B7").Value = "'5/6"
$ws.Range("C7").Value = "phương nguyên"
